# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append new game's yardage-per-play log to the running totals ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = '3 3 1 1 18 7 9 9 6 8 8 4 4 3 -1 11 2 27 3 4 3 2 9 0 -2 4 -1 4 14 -1 2 0 17 5 0 14 3 15 10 4 1 -1 4 3 1 21 1 5 2 4 2 3 22 1 3 4 5 1 1 -1 9 4 16 1 7 3 11 2 0 -1 0 12 6 5 3 3 1 2 4 2 2 5 6 5 1 3 13 1 9 4 2 1 5 -1 11 4 1 2 1 18 9 4 2 8 5 6 6 3 3 5 10 8 5 1 4 6 2 -1 2 3 0 1 5 3 -1 9 1 9 4 2 31 0 5 5 5 12 14 5 11 5 3 4 6 8 3 17 0 8 11 5 1 8 4 5 5 7 3 9 0 -3 2 -4 3 0 9 -4 2 13 2 4 11 16 1 2 27 0 1 16 3 1 1 0 0 2 3 1 3 2 4 1 4 -3 4 3 4 3 3 8 1 2 0 11 -2 2 2 0 2 7 4 8 0 3 0 4 5 3 8 6 10 4 3 5 11 3 2 9 2 3 0 -2 14 -4 -5 7 2 3 6 1 0 8 -2 2 3 4 10 10 17 5 0 0 5 4 7 4 4 8 1 1 4 2 5 6 30 20 3 0 16 3 1 6 3 10 6 5 12 4 1 -3 2 6 3 13 32 0 8 5 9 -2 3 -1 9 2 12 4 0 0 0 -3 -8 1 5 1 1 13 -5 3 8 4 11 6 7 1 6 3 16 1 3 4 1 2 1 6 3 2 1 12 4 0 8 24 4 4 1 2 5 5 10 -3 0 3 1 13 9 7 3 6 20 11 4 3 11 2 1 17 4 3 3 2 3 6 3 0 1 0 10 3 2 6 5 2 3 1 4 12 6 4 2 1 2 4 9 13 4 1 0 16 4 1 9 2 2 13 5 -2 50 6 6 2 4 1 -3 3 9 0 5 11 1 2 5 2 -3 13 11 2 5 5 11 6 -1 2 10 7 26 10 3 3 7 1 -2 4 4 1 2 5 3 4 4 3 8 9 7 1 2 5 4 7 0 4 0 9 1 7 -1 1 3 9 4 2 7 6 15 1 3 -5 4 -2 -2 4 12 10 4 11 7 8 7 3 4 7 5 14 10 4 2 3 8 4 2 18 2 12 7 5 7 4 -1 2 2 21 5 4 5 4 5 3 12 2 15 24 17 1 1 2 9 3 10 10 1 8 1 2 4 1 12 9 16 3 3 5 4 4 2 10 5 2 3 15 3 6 4 1 9 5 0 -1 23 9 -1 -2 9 5 9 0 7 2 0 6 3 3 17 3 6 1 6 4 0 5 11 4 3 3 3 10 5 2 -1 1 3 0 5 2 2 11 -1 13 18 6 9 9 6 2 5 2 10 -3 4 2 4 4 9 7 5 9 3 8 12 -1 3 4 6 0 -1 1 8 6 -1 7 7 2 6 0 0 6 4 7 7 7 2 4 10 2 -4 -1 -1 2 3 2 5 0 5 1 3 21 1 3 2 0 6 -1 5 3 2 6 -1 4 4 4 3 6 8 33 11 4 3 4 1 5 1 5 0 2 3 3 7 4 14 3 6 4 10 -2 -7 10 2 3 5 10 4 2 12 4 2 0 5 4 2 4 4 4 4 6 7 -4 1 2 1 2 7 5 7 0 3 1 0 0 1 4 9 6 4 2 10 4 2 -2 2 4 2 5 2 5 11 9 7 0 5 0 -1 3 3 1 2 3 0 1 5 0 1 32 0 3 4 0 1 3 8 12 5 5 2 3 2 7 2 6 3 8 3 10 4 6 7 2 7 4 0 1 0 0 1 5 -2 3 1 2 23 -1 6 2 1 8 10 6 8 1 14 3 17 1 19 8 0 5 8 11 0 11 3 3 3 11 3 0 1 6 10 5 6 5 25 5 3 6 7 4 6 4 4 2 5 6 4 4 4 2 2 -1 23 -2 15 13 12 1 0 -2 -4 5 4 3 6 11 -1'
$ydsWs.Range("C2").Value = '5 7 9 2 19 4 6 1 4 7 2 2 2 1 5 6 17 2 4 13 1 9 9 6 4 3 5 3 0 1 3 3 5 3 13 5 3 1 7 9 5 8 7 2 3 3 5 9 1 8 7 4 2 0 7 3 3 -1 1 11 1 -1 2 5 3 30 5 6 11 3 6 4 18 4 8 11 1 4 5 24 3 2 8 5 2 2 3 4 10 18 -1 2 1 2 3 6 2 5 6 -1 8 4 4 7 0 9 41 4 2 5 3 5 1 10 7 3 6 1 2 7 2 2 43 1 3 4 4 2 0 4 2 2 0 7 7 2 4 7 2 5 1 9 10 3 9 4 2 3 5 4 1 13 3 -3 3 6 -1 7 -6 2 0 6 3 -1 6 3 2 14 2 6 20 9 4 19 6 1 1 6 5 10 8 -1 10 6 4 2 4 11 0 3 2 3 3 3 6 7 6 4 0 3 9 3 5 2 2 3 1 6 4 2 5 0 5 6 0 8 5 2 6 5 4 7 5 3 6 13 4 6 0 10 3 -1 1 2 -4 4 4 1 6 1 1 4 13 2 2 2 5 6 5 3 0 2 5 2 2 11 1 0 4 2 5 12 4 0 -3 5 4 1 2 11 3 1 34 1 6 0 3 1 0 2 8 0 4 8 8 4 3 65 1 6 3 3 13 2 0 3 5 2 7 1 3 -1 4 5 1 8 6 1 1 7 7 1 -2 6 2 2 2 2 2 3 5 2 3 7 2 9 3 5 2 1 9 5 5 4 2 3 7 2 1 11 2 3 6 2 4 5 2 1 8 14 4 3 4 -7 5 0 6 3 4 6 6 3 3 2 2 4 3 2 6 5 2 1 3 1 9 2 5 2 0 1 11 0 8 5 3 1 9 3 3 34 1 6 2 -2 5 5 5 6 0 6 4 -3 23 18 2 9 3 7 6 1 10 4 4 3 0 -1 1 3 7 17 2 15 3 7 10 4 4 5 1 18 8 11 3 13 -3 3 6 6 11 7 8 4 1 1 0 3 4 27 6 3 2 0 4 7 8 4 4 8 1 1 5 2 6 7 -1 6 4 2 4 17 8 18 5 15 6 17 4 2 6 7 6 2 2 -2 5 2 6 5 12 31 9 7 1 9 9 8 0 0 21 3 4 5 20 14 -4 3 4 -2 20 6 2 8 5 9 2 1 11 5 9 -1 3 1 1 2 1 4 -1 9 1 2 3 3 5 -2 -1 5 9 1 -1 11 6 9 6 0 8 1 2 5 11 12 5 12 6 3 3 -4 7 2 2 7 4 2 13 10 0 8 22 2 11 9 9 8 2 -1 3 5 -1 0 0 3 1 4 12 4 6 6 0 8 3 4 8 1 3 5 3 3 2 2 15 1 7 8 6 4 4 11 2 5 3 9 2 3 4 0 8 5 0 3 4 -1 2 4 0 11 1 -1 2 2 4 6 2 1 2 8 5 4 -4 9 3 2 -4 5 2 0 2 -3 11 2 2 6 2 14 5 4 0 3 5 4 10 2 2 4 9 15 2 3 -1 7 7 -2 8 5 3 1 5 5 0 7 8 2 2 4 9 6 4 9 3 5 3 3 1 -1 12 2 -1 9 2 1 3 8 -1 2 5 31 3 1 5 1 5 2 8 3 5 4 4 4 9 -1 1 14 2 9 6 7 2 14 3 0 1 11 3 5 5 -1 5 4 8 9 7 2 6 2 16 1 -3 5 1 4 3 4 1 5 3 4 15 3 15 4 3 6 2 4 1 2 7 0 9 4 6 8 4 3 1 8 2 4 8 22 1 1 3 0 2 2 3 2 6 15 2 2 -2 0 2 1 21 5 5 8 6 6 6 12 4 6 3 6 6 5 8 4 2 11 17 2 13 3 6 6 3 1 4 3 3 1 3 4 5 4 1 0 0 0 6 2 6 1 9 6 1 5 29 6 5 3 0 3 23 -1 4 3 -3 47 7 0 5 15 8 -4 4 1 2 -2 2 5 2 4 3 4 4 8 0 -2 2 4 7 2 13 4'
$ydsWs.Range("B3").Value = '6 6 11 7 15 10 13 2 6 5 4 9 9 2 17 11 7 12 14 7 19 15 1 3 13 1 12 5 10 6 7 20 9 14 13 10 11 22 2 54 11 16 8 14 8 4 5 4 11 33 15 11 13 15 2 5 7 29 5 18 -4 20 10 49 9 5 7 24 7 6 9 7 6 18 13 11 2 17 5 8 22 17 11 2 14 21 24 8 5 5 10 6 19 45 6 6 6 1 9 12 26 32 4 13 5 16 37 8 4 23 15 23 6 16 2 42 4 3 -1 15 7 8 5 6 22 11 5 15 4 12 15 16 8 9 6 -3 21 -2 20 2 37 8 11 8 19 14 38 4 5 15 17 7 6 17 7 22 10 11 2 13 9 15 18 30 5 13 2 36 7 4 15 23 -3 18 5 3 -4 14 9 25 14 16 8 6 8 19 26 16 5 41 10 7 0 14 10 3 8 6 32 23 7 5 8 14 1 10 9 10 28 14 29 -5 3 44 4 10 26 28 8 13 2 3 8 10 14 6 8 4 3 17 21 11 8 3 19 11 6 3 9 8 11 12 13 3 5 15 6 19 9 10 9 16 15 3 22 34 23 6 3 1 75 19 3 14 7 44 20 12 14 2 9 5 13 9 2 8 29 37 7 7 8 9 14 8 2 2 4 6 20 8 11 3 15 3 2 8 19 14 18 20 1 19 6 4 10 -5 15 28 5 20 12 15 30 14 5 20 8 5 17 9 9 1 16 17 7 37 21 9 13 16 9 5 15 6 26 44 28 26 5 32 13 5 11 9 8 4 -5 11 6 5 14 10 23 23 19 1 2 13 4 20 9 5 16 9 9 8 16 7 7 10 15 6 10 31 5 4 13 9 9 5 11 2 9 4 11 12 13 26 0 4 31 2 9 17 5 19 16 9 25 12 21 9 -2 37 3 7 2 3 11 8 8 18 12 3 -2 6 9 12 10 8 10 12 31 8 3 7 26 -1 9 42 9 11 20 8 14 8 15 -5 9 7 3 13 22 16 14 8 5 23 24 4 5 1 9 1 11 11 16 11 8 3 3 5 15 9 11 33 17 16 5 6 8 9 7 0 71 1 12 12 5 3 6 14 13 4 12 11 -1 5 8 0 5 21 13 17 33 12 9 4 23 8 18 16 1 4 17 9 11 24 9 8 -6 4 11 6 22 19 -2 9 28 6 7 12 6 11 8 3 14 75 8 4 14 6 7 19 10 33 17 11 14 10 4 14 20 14 5 11 40 5 46 3 8 3 23 13 7 7 15 5 6 6 8 17 -1 5 15 28 7 6 11 12 8 10 2 6 23 18 10 5 5 7 11 8 11 9 36 1 9 2 19 32 4 6 4 9 2 15 7 6 3 14 7 10 5 12 11 44 12 6 15 11 9 7 8 5 8 10 6 26 8 12 8 11 5 4 3 11 8 17 6 8 5 5 4 1 8 6 4 4 11 6 20 9 19 27 49 3 5 7 27 19 17 14 9 10 5 13 9 5 12 13 2 9 0 3 11 5 21 19 6 12 6 24 7 8 6 6 5 10 9 23 8 11 24 12 14 6 4 9 10 4 2 10 3 9 12 7 0 15 5 6 6 7 11 9 17 9 3 11 12 9 18 8 4 9 10 10 15 14 24 2 2 7 10 4 2 1 7 2 7 12 2 25 6 14 10 13 22 9 13 6 7 5 11 22 19 8 8 7 13 12 -1 16 27 3 1 0 9 7 11 7 13 17 2 10 16 32 3 7 3 38 8 30 22 1 10 8 12 0 8 37 13 23 7 11 8 24 14 22 10 12 6 6 4 8 8 8 9 9 12 13 9 29 14 -7 12 6 14 14 38 8 4 6 44 23 10 19 19 28 7 9 3 1 7 38 15 5 -3 2 8 16 1 12 14 5 9 20 11 12 -1 17 4 3 6 5 23 10 40 5 3 9 0 6 69 1 5 10 20 7 10 27 14 34 10 11 11 9 13 -1 5 18 18 5 7 3 14 8 50 3 8 9 4 6 9 16 22 8 6 1 29 7 17 24 3 3 53 7 3 2 4 27 7 8 5 11 6 5 6 4 8 1 4 3 21 11 4 9 4 3 3 11 11 10 11 2 11 17 10 7 14 5 10 4 1 19 4 44 8 12 -4 20 5 5 4 15 20 2 14 15 5 4 31 0 3 12 15 27 48 29 8 23 4 1 31 6 15 3 41 2 2'
$ydsWs.Range("C3").Value = '12 5 6 11 9 19 7 15 20 8 7 20 10 12 19 31 17 8 6 35 16 14 10 12 5 22 6 14 11 14 8 32 17 14 25 22 7 5 19 1 6 7 0 0 7 4 6 13 7 7 8 5 19 6 6 4 7 6 13 11 0 8 8 9 15 -4 5 8 19 0 25 4 12 30 9 46 12 5 59 23 5 72 3 8 10 1 4 2 4 23 7 2 4 11 6 4 2 15 13 4 7 5 14 22 8 3 27 5 2 10 7 19 5 17 13 21 12 13 4 6 6 13 11 24 5 11 9 8 4 15 5 -2 1 27 13 5 15 8 1 5 -6 10 6 8 3 7 4 1 12 7 0 5 11 7 14 9 3 11 12 -6 8 8 3 28 24 3 11 2 9 22 11 5 3 9 6 23 0 26 29 16 7 11 17 18 14 11 11 4 1 18 16 7 9 3 21 3 20 5 7 1 12 6 6 2 6 8 48 -4 44 9 14 9 8 7 5 7 8 6 10 11 7 37 7 5 8 11 0 10 5 5 8 21 11 10 8 5 16 14 18 7 7 10 14 2 12 16 9 10 10 -1 11 11 13 5 23 10 29 2 6 9 9 18 4 22 5 51 21 8 25 24 9 14 8 12 5 12 17 26 -1 3 15 7 7 5 8 54 5 7 0 15 17 6 9 6 2 3 7 7 17 11 8 22 9 4 20 11 32 6 6 12 17 5 8 48 27 18 27 11 -3 3 7 -8 26 26 23 25 18 0 4 2 6 17 4 3 4 5 2 17 6 7 3 2 6 14 12 20 7 2 6 17 23 6 8 6 10 6 13 9 11 7 15 6 6 34 1 5 16 15 5 8 3 31 -3 9 4 17 8 5 15 1 3 12 25 15 6 8 14 16 8 44 17 3 9 10 9 43 11 9 30 4 10 3 4 19 11 17 -4 7 20 8 9 7 11 4 4 12 8 27 20 42 19 9 10 6 0 4 3 8 8 11 19 4 20 8 6 16 5 12 7 21 11 5 0 20 43 15 6 16 4 13 22 21 12 10 6 3 5 15 3 7 14 18 6 37 12 1 21 8 7 6 4 14 14 9 18 19 17 9 9 15 23 8 24 35 61 15 53 41 1 -2 11 5 16 16 8 9 13 6 3 7 12 2 7 3 5 6 6 6 4 16 39 10 14 4 -4 -4 9 12 4 -2 7 10 18 3 5 3 24 8 13 22 10 2 20 6 12 9 11 6 5 8 24 5 50 19 3 11 1 12 41 2 5 4 22 4 6 11 6 5 7 16 4 3 3 10 15 5 6 21 35 0 11 -1 14 16 3 8 2 15 5 6 5 1 4 3 27 3 6 7 5 22 37 24 40 11 5 5 7 9 6 12 3 15 -2 8 -3 11 11 2 2 8 11 3 1 5 4 5 16 4 0 6 12 18 14 -2 3 19 5 19 7 11 16 36 10 10 6 3 2 34 5 9 5 8 11 9 13 19 12 6 19 10 7 7 8 9 16 14 2 7 19 7 2 5 7 13 10 6 4 6 4 19 19 6 4 9 4 12 12 13 15 4 2 27 7 5 4 17 6 29 8 5 15 15 8 5 14 18 9 5 -3 24 10 0 1 3 5 10 6 2 17 4 8 2 5 6 6 3 3 7 4 6 12 11 4 3 10 72 8 8 9 17 18 6 7 6 4 10 5 69 22 6 4 39 10 5 2 5 9 2 35 30 6 16 11 8 2 31 11 6 13 7 10 29 28 3 4 7 5 5 0 16 0 6 8 5 5 22 13 1 4 7 9 8 15 5 5 7 7 10 12 11 11'

# --- OFF sheet: season offensive totals updated with the new game ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value = 5
$offWs.Range("C2").Value = 456
$offWs.Range("D2").Value = 26
$offWs.Range("E2").Value = 16
$offWs.Range("F2").Value = 100
$offWs.Range("G2").Value = 128
$offWs.Range("I2").Value = 15
$offWs.Range("J2").Value = 59
$offWs.Range("L2").Value = 831
$offWs.Range("M2").Value = 543
$offWs.Range("O2").Value = 36
$offWs.Range("P2").Value = 22
$offWs.Range("Q2").Value = 1369

$offWs.Range("C3").Value = 499
$offWs.Range("D3").Value = 21
$offWs.Range("E3").Value = 72
$offWs.Range("F3").Value = 285
$offWs.Range("G3").Value = 137
$offWs.Range("H3").Value = 59
$offWs.Range("I3").Value = 125
$offWs.Range("J3").Value = 130
$offWs.Range("N3").Value = 33

# --- DEF sheet: season defensive totals updated with the new game ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 425
$defWs.Range("F2").Value = 144
$defWs.Range("G2").Value = 128
$defWs.Range("I2").Value = 13
$defWs.Range("J2").Value = 73
$defWs.Range("L2").Value = 713
$defWs.Range("M2").Value = 459
$defWs.Range("O2").Value = 60
$defWs.Range("P2").Value = 39
$defWs.Range("Q2").Value = 1268

$defWs.Range("B3").Value = 21
$defWs.Range("C3").Value = 422
$defWs.Range("E3").Value = 64
$defWs.Range("F3").Value = 234
$defWs.Range("G3").Value = 114
$defWs.Range("H3").Value = 48
$defWs.Range("I3").Value = 118
$defWs.Range("J3").Value = 130
$defWs.Range("N3").Value = 42

# --- ST sheet: special teams totals + logged per-play kick/return/penalty-yard distributions ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 216
$stWs.Range("D2").Value = 97
$stWs.Range("F2").Value = 257
$stWs.Range("G2").Value = 241
$stWs.Range("H2").Value = 6
$stWs.Range("I2").Value = 5
$stWs.Range("J2").Value = 95
$stWs.Range("K2").Value = 92
$stWs.Range("N2").Value = 28
$stWs.Range("O2").Value = 20
$stWs.Range("B3").Value = 153

$stWs.Range("B4").Value = '66 66 63 60 64 69 63 57 58 65 65 51 44 68 65 62 67 63 46 66 61 61 46 60 67 68 64 60 59 73 60 63 65 12 55 62 61 66 66 60 66 67 62 46 63 64 65 63 60 60 66 63 68 58 63 64 64 64 63 32 63 50 64'
$stWs.Range("D3").Value = '44 38 54 38 42 32 54 55 58 65 65 53 60 44 67 53 47 33 33 42 38 57 29 38 44 40 36 55 52 50 54 32 23 39 30 43 41 58 48 51 34 61 45 37 41 25 31 52 58 24 24 44 38 44 51 27 29 44 47 54 49 32 51 46 42 34 19 52 41 50 52 43 54 55 44 64 42 49 49 68 31 54 50 37 52 32 41 50 42 35 31 58 43 41 33 32 32'
$stWs.Range("D4").Value = '0 0 19 0 3 0 0 0 16 12 0 0 0 0 0 0 0 0 0 24 0 0 0 0 0 0 0 8 3 0 2 0 0 -2 0 0 0 -4 11 1 0 17 9 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 16 0 0 7 0 0 0 0 0 0 0 0 0 0 0 15 0 0 7 0 0 14 6 0 6 0 0 20 0 3 0 0 0 8 0 0 5'
$stWs.Range("B5").Value = '21 22 23 16 28 20 29 26 93 19 19 18 9 29 17 23 23 0 6 22 31 24 7 18 27 26 28 21 13 17 20 14 44 0 15 23 19 17 17 12 19 26 19 3 23 19 18 19 21 18 26 18 28 16 29 75 19 29 13 0 23 14 29'
$stWs.Range("D5").Value = '0 0 0 0 11 8 18 0 0 29 3 5 0 0 21 1 48 0 0 0 -8'
$stWs.Range("B6").Value = '20 31 19 22 23 22 11 2 23 26 29 21 29 27 15 21 44 26 19 31 28 18 24 18 26 41 21 25 24 20 8 33 3 9 25 30 28 22 32 28 31 33 0 28 23 12 20 14 38 27 39 24 19 20 25 26 30 12 23 27'

# --- TURNS sheet: turnover totals updated with the new game ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value = 14
$turnsWs.Range("C2").Value = 19
$turnsWs.Range("D2").Value = 15
$turnsWs.Range("E2").Value = 21
$turnsWs.Range("E3").Value = 19

# --- PEN sheet: penalty totals updated with the new game ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 43
$penWs.Range("B3").Value = 42
